$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the now-unneeded "Shared uncertainty" duration/correlation data for rows 4-7
$ws.Range("X4:AC7").ClearContents()

# Scroll/select as the user left it after performing the edit
$ws.Range("X4:AC7").Select()
$excel.ActiveWindow.ScrollColumn = 18
